$d = $word.ActiveDocument

function New-RunXml([string]$innerBody) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Step 1: split the first "To take advantage..." paragraph (paragraph 3) into
# two runs, breaking right after "like Emphasis and ".
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$run1a = "To take advantage of this template’s design, use the Styles gallery on the Home tab. You can format your headings by using heading styles, or highlight important text using other styles, like Emphasis and "
$run1b = "Intense Quote. These styles come in formatted to look great and work together to help communicate your ideas."
$xml1 = New-RunXml("<w:p><w:r><w:t xml:space=`"preserve`">$run1a</w:t></w:r><w:r><w:t>$run1b</w:t></w:r></w:p>")
$r3.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Step 2: merge the "Text" / " Signatures" runs (paragraph 6 heading) into one
# run and prefix it with a lastRenderedPageBreak marker.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$r6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$xml2 = New-RunXml("<w:p><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Text Signatures</w:t></w:r></w:p>")
$r6.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Step 3: split the second "To take advantage..." paragraph (paragraph 7)
# into two runs, breaking after "format yo".
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$r7 = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$run2a = "To take advantage of this template’s design, use the Styles gallery on the Home tab. You can format yo"
$run2b = "ur headings by using heading styles, or highlight important text using other styles, like Emphasis and Intense Quote. These styles come in formatted to look great and work together to help communicate your ideas."
$xml3 = New-RunXml("<w:p><w:r><w:t>$run2a</w:t></w:r><w:r><w:t>$run2b</w:t></w:r></w:p>")
$r7.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Step 4: remove the second page break, the third "Text Signatures" heading
# and the third "To take advantage..." paragraph (paragraphs 9-11), merging
# their removal into paragraph 12 (the one holding the bookmark).
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$p12 = $d.Paragraphs.Item(12)
$rDel1 = $d.Range($p9.Range.Start, $p12.Range.Start)
$rDel1.Cut()

# ---------------------------------------------------------------------------
# Step 5: remove the "Go ahead and get started." text that now precedes the
# _GoBack bookmark, leaving the bookmark markers in place.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$pBookmark = $d.Paragraphs.Item(9)
$rDel2 = $d.Range($pBookmark.Range.Start, $bm.Start)
$rDel2.Cut()

# ---------------------------------------------------------------------------
# Step 6: drop the trailing empty paragraph. Delete the (now empty besides
# the bookmark) paragraph together with its mark, then restore the bookmark
# inside what becomes the new final paragraph of the document.
# ---------------------------------------------------------------------------
$pBookmark2 = $d.Paragraphs.Item(9)
$pTrailing = $d.Paragraphs.Item(10)
$rDel3 = $d.Range($pBookmark2.Range.Start, $pTrailing.Range.Start)
$rDel3.Cut()

$pFinal = $d.Paragraphs.Item($d.Paragraphs.Count)
$rFinal = $d.Range($pFinal.Range.Start, $pFinal.Range.Start)
$xml4 = New-RunXml("<w:p><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>")
$rFinal.InsertXML($xml4)
